# Update the "prediction" values in column B (rows 2-25) with the refreshed
# scores from the regenerated outputs-r202 run (ful-path.csv re-export).
# Column A (Row labels) and column C (unchanged 1's) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1.2591089106219577
    3  = 0.76958862334451261
    4  = 0.63064843470889187
    5  = 1.2130288751369047
    6  = 0.58232450978965744
    7  = 1.1416662459255598
    8  = 1.0211669964976497
    9  = 0.76351749550507364
    10 = 0.74474099773206426
    11 = 0.97225023792195664
    12 = 0.72719953982875118
    13 = 1.1780946490355273
    14 = 0.78927200538351805
    15 = 0.68524069468210413
    16 = 0.83150974979921521
    17 = 0.66038858707405979
    18 = 1.1349131961421559
    19 = 1.2843023351137663
    20 = 1.2814771998226497
    21 = 1.308379608895228
    22 = 0.71540685255352443
    23 = -0.15431479402118065
    24 = -0.55197715686269611
    25 = 0.52712928886984933
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
